# Weekly update: insert a new price record as the new top row (row 193) for
# the "Ají" (Hortaliza) block on Sheet1, shifting the existing rows 193:215
# down to 194:216.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 193; everything below (old rows 193:215) moves
# down to 194:216, carrying its original data/format with it.
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new weekly record.
$ws.Range("A193").Value = 8
$ws.Range("B193").Value = "Terminal La Palmera de La Serena"
$ws.Range("C193").Value = "Coquimbo"
$ws.Range("D193").Value = 44694
$ws.Range("E193").Value = 4
$ws.Range("F193").Value = 100112021
$ws.Range("G193").Value = "Ají"
$ws.Range("H193").Value = "Americana (o)"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 200
$ws.Range("K193").Value = 27000
$ws.Range("L193").Value = 28000
$ws.Range("M193").Value = 27500
$ws.Range("N193").Value = "`$/caja 25 kilos"
$ws.Range("O193").Value = "Provincia de Limarí"
$ws.Range("P193").Value = 1100
$ws.Range("Q193").Value = 25
$ws.Range("R193").Value = "Hortaliza"

# D193 carries a date style; match the style used by the other date cells
# ("YYYY-MM-DD HH:MM:SS", the same numeric format as every other Fecha cell).
$ws.Range("D193").NumberFormat = "YYYY-MM-DD HH:MM:SS"
